$d = $word.ActiveDocument

# --- Paragraph 3: Financial status ---
$p3 = $d.Paragraphs(3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$r3.Text = "Financial status: I have paid for all my parts needed, I have an extra sound sensor just in case, my last purchases were a cheap USB keyboard and USB mouse, which adds 10 dollars to my expenses, I needed these purchases to control the raspberry PI."

# --- Paragraph 4: Blog/Documentation ---
$p4 = $d.Paragraphs(4)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$r4.Text = "Blog/Documentation: My blog is up to date, I documented all my activities there. "

# --- Paragraph 5: was Milestone, becomes the blog link paragraph ---
$p5 = $d.Paragraphs(5)
$r5 = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$r5.Text = "This is the link to my blog: https://github.com/SlavaPere/SensorEffector/blob/master/index.md"

# --- Paragraph 6: was Activities, becomes the Milestone paragraph ---
$p6 = $d.Paragraphs(6)
$r6 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$r6.Text = "Milestone: I" + [char]0x2019 + "m ready to demonstrate the basic capabilities of my sound sensor, I can graphically display different levels of sound captured on the screen. I have met prior milestones."

# --- Paragraph 7: was Problems, becomes the Activities paragraph (keeps _GoBack bookmark) ---
$p7 = $d.Paragraphs(7)
$p7Start = $p7.Range.Start
$r7 = $d.Range($p7Start, $p7.Range.End - 1)
$beforeBookmark = "Activities: I" + [char]0x2019 + "m starting work on Fritzing"
$afterBookmark = "; I" + [char]0x2019 + "m planning to add a statement in the code that can detect a persistent noise level (sound over a threshold level over a period of time). I" + [char]0x2019 + "m getting ready to add a second sound sensor as well that will detect a different pattern, I also am aware that I have to write the build instructions soon."
$r7.Text = $beforeBookmark + $afterBookmark
$bmPos = $p7Start + $beforeBookmark.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Insert a new paragraph 8 after paragraph 7: the new Problems paragraph ---
$p7 = $d.Paragraphs(7)
$insertPoint = $d.Range($p7.Range.End - 1, $p7.Range.End - 1)
$insertPoint.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)
$r8 = $d.Range($p8.Range.Start, $p8.Range.End - 1)
$problemsText = "Problems: As of 12/21/2017 the most pressing problem I have is that my sound sensor doesn" + [char]0x2019 + "t seem to work with an integrated PCB, although the sound sensor and PCF get powered the sound level from my program is constantly displaying " + [char]0x201C + "-1" + [char]0x201D + " which doesn" + [char]0x2019 + "t make sense (It means the mike is reading too high levels, as the lower the number - the higher the decibel level it implies. I first plan to check if my PCB lost some connectivity, if it didn" + [char]0x2019 + "t, I will have to look into a calibration, wiring or coding/support issue."
$r8.Text = $problemsText
